# Apply weekly re-shuffle of Fecha/Volumen/Precio rows (rows 2-9) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
$rows = @(2, 3, 4, 5, 6, 7, 8, 9)
$D = @{ 2 = 44497; 3 = 44475; 4 = 44489; 5 = 44461; 6 = 44455; 7 = 44454; 8 = 44490; 9 = 44482 }
$M = @{ 2 = 500;   3 = 240;   4 = 160;   5 = 200;   6 = 200;   7 = 160;   8 = 400;   9 = 240 }
$N = @{ 2 = 9000;  3 = 11000; 4 = 9500;  5 = 11000; 6 = 12000; 7 = 12000; 8 = 9500;  9 = 10000 }
$O = @{ 2 = 10000; 3 = 12000; 4 = 10000; 5 = 12000; 6 = 13000; 7 = 13000; 8 = 10000; 9 = 11000 }
$P = @{ 2 = 9500;  3 = 11500; 4 = 9750;  5 = 11500; 6 = 12500; 7 = 12500; 8 = 9750;  9 = 10500 }
$S = @{ 2 = 4750;  3 = 5750;  4 = 4875;  5 = 5750;  6 = 6250;  7 = 6250;  8 = 4875;  9 = 5250 }

foreach ($r in $rows) {
    $ws.Range("D$r").Value = $D[$r]
    $ws.Range("M$r").Value = $M[$r]
    $ws.Range("N$r").Value = $N[$r]
    $ws.Range("O$r").Value = $O[$r]
    $ws.Range("P$r").Value = $P[$r]
    $ws.Range("S$r").Value = $S[$r]
}
